$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Q8 in J1 (headers already exist B1:I1 = Q0..Q7)
# Copy formatting from I1 (bordered/bold/centered header style) onto J1
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "Q8"

# Row 2 - new values (columns B:H, same extent as before)
$ws.Range("B2").Value = -0.2634544738301702
$ws.Range("C2").Value = -0.1677772249736381
$ws.Range("D2").Value = -0.4281992330346231
$ws.Range("E2").Value = 0.4072587527989196
$ws.Range("F2").Value = 0.3072359837017114
$ws.Range("G2").Value = -0.04662207633221729
$ws.Range("H2").Value = 0.1015709931031478

# Row 3
$ws.Range("B3").Value = -0.1427655701645288
$ws.Range("C3").Value = -0.4031875782255139
$ws.Range("D3").Value = 0.4322704076080288
$ws.Range("E3").Value = 0.3322476385108206
$ws.Range("F3").Value = -0.02161042152310803
$ws.Range("G3").Value = 0.1265826479122571

# Row 4 - grows to J4
$ws.Range("B4").Value = -0.1674069739471072
$ws.Range("C4").Value = 0.6680510118864355
$ws.Range("D4").Value = 0.5680282427892273
$ws.Range("E4").Value = 0.2141701827552986
$ws.Range("F4").Value = 0.3623632521906637
$ws.Range("G4").Value = 0.519796913698928
$ws.Range("H4").Value = 0.0191408274226057
$ws.Range("I4").Value = 0.2994971146932124
$ws.Range("J4").Value = 0.1275310031338272

# Row 5
$ws.Range("B5").Value = 1.641111674953823
$ws.Range("C5").Value = 1.541088905856615
$ws.Range("D5").Value = 1.187230845822686
$ws.Range("E5").Value = 1.335423915258051
$ws.Range("F5").Value = 1.492857576766315
$ws.Range("G5").Value = 0.9922014904899932
$ws.Range("H5").Value = 1.2725577777606
$ws.Range("I5").Value = 1.100591666201215

# Row 6
$ws.Range("B6").Value = 0.64310306534761
$ws.Range("C6").Value = 0.2892450053136812
$ws.Range("D6").Value = 0.4374380747490463
$ws.Range("E6").Value = 0.5948717362573106
$ws.Range("F6").Value = 0.09421564998098833
$ws.Range("G6").Value = 0.374571937251595
$ws.Range("H6").Value = 0.2026058256922098

# Row 7
$ws.Range("B7").Value = 0.2251646141706316
$ws.Range("C7").Value = 0.3733576836059967
$ws.Range("D7").Value = 0.530791345114261
$ws.Range("E7").Value = 0.03013525883793872
$ws.Range("F7").Value = 0.3104915461085453
$ws.Range("G7").Value = 0.1385254345491602

# Row 8 - grows to I8
$ws.Range("B8").Value = 0.5378291618471884
$ws.Range("C8").Value = 0.6952628233554528
$ws.Range("D8").Value = 0.1946067370791305
$ws.Range("E8").Value = 0.4749630243497371
$ws.Range("F8").Value = 0.302996912790352
$ws.Range("G8").Value = 0.6087440311202421
$ws.Range("H8").Value = -0.003446134512213206
$ws.Range("I8").Value = 0.6848900978272674

# Row 9
$ws.Range("B9").Value = 0.5086525859949649
$ws.Range("C9").Value = 0.007996499718642608
$ws.Range("D9").Value = 0.2883527869892493
$ws.Range("E9").Value = 0.1163866754298641
$ws.Range("F9").Value = 0.4221337937597542
$ws.Range("G9").Value = -0.1900563718727011
$ws.Range("H9").Value = 0.4982798604667795

# Row 10
$ws.Range("B10").Value = -0.1889757213638347
$ws.Range("C10").Value = 0.09138056590677199
$ws.Range("D10").Value = -0.08058554565261317
$ws.Range("E10").Value = 0.2251615726772769
$ws.Range("F10").Value = -0.3870285929551783
$ws.Range("G10").Value = 0.3013076393843023

# Row 11
$ws.Range("B11").Value = 0.1816619472563389
$ws.Range("C11").Value = 0.009695835696953689
$ws.Range("D11").Value = 0.3154429540268438
$ws.Range("E11").Value = -0.2967472116056115
$ws.Range("F11").Value = 0.3915890207338691

# Row 12
$ws.Range("B12").Value = -0.1282989989530805
$ws.Range("C12").Value = 0.1774481193768096
$ws.Range("D12").Value = -0.4347420462556457
$ws.Range("E12").Value = 0.2535941860838349

# Row 13
$ws.Range("B13").Value = 0.2042627890992136
$ws.Range("C13").Value = -0.4079273765332416
$ws.Range("D13").Value = 0.280408855806239

# Row 14
$ws.Range("B14").Value = -0.4731698452888153
$ws.Range("C14").Value = 0.2151663870506653

# Row 15
$ws.Range("B15").Value = 0.2908260759093906

# Row 16 - unchanged (only A16 label, no numeric values)
